$d = $word.ActiveDocument

$replacements = @(
    @("936×7=6552", "111×4=444"),
    @("184×8=1472", "300×4=1200"),
    @("864×3=2592", "633×2=1266"),
    @("808×9=7272", "900×8=7200"),
    @("681×5=3405", "441×2=882"),
    @("382×2=764", "710×8=5680"),
    @("181×2=362", "635×6=3810"),
    @("337×6=2022", "407×8=3256"),
    @("401×7=2807", "807×4=3228"),
    @("549×2=1098", "734×4=2936"),
    @("657×7=4599", "443×6=2658"),
    @("760×7=5320", "687×8=5496"),
    @("446×9=4014", "977×6=5862"),
    @("933×6=5598", "817×8=6536"),
    @("357×7=2499", "397×5=1985"),
    @("931×2=1862", "484×8=3872"),
    @("385×6=2310", "333×9=2997"),
    @("766×5=3830", "739×4=2956"),
    @("477×7=3339", "340×5=1700"),
    @("287×4=1148", "471×7=3297"),
    @("943×4=3772", "455×5=2275"),
    @("974×7=6818", "176×6=1056"),
    @("687×9=6183", "362×3=1086"),
    @("446×3=1338", "416×7=2912"),
    @("830×6=4980", "627×6=3762")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
